$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.602.56'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.672.09'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.35%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.22%  '
$ws.Range('E6').Value = '  +1.84%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.17'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.263'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0638'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0903'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.916.18'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.670.49'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.603'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.85'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.00'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.652.90'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.07'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.30'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0718'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.20%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.22'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.92'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.112'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.76'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.16%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0492'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.45'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.32'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.494.91'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.76'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.25%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.02'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('B37').Value = 'Aave'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '83.11'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.595'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.81%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0177'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.18%  '
$ws.Range('E40').Value = '  -3.27%  '
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.834'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0499'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('E45').Value = '  +1.89%  '
$ws.Range('E46').Value = '  +0.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.53'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.809.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '49.77'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '93.26'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0112'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.57%  '
